$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.184.81"
$ws.Range("E2").Value = "  -3.91%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.492.29"
$ws.Range("E3").Value = "  -5.39%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.14"
$ws.Range("E5").Value = "  -6.74%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.68"
$ws.Range("E6").Value = "  -8.30%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.490.48"
$ws.Range("E7").Value = "  -5.40%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  -3.58%  "

# Row 10
$ws.Range("E10").Value = "  -5.87%  "

# Row 11
$ws.Range("E11").Value = "  -3.95%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.422"
$ws.Range("E12").Value = "  -5.08%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000215"
$ws.Range("E13").Value = "  -7.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.081.48"
$ws.Range("E14").Value = "  -5.37%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.43"
$ws.Range("E15").Value = "  -4.31%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.478.62"
$ws.Range("E16").Value = "  -5.71%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.076.49"
$ws.Range("E17").Value = "  -4.04%  "

# Row 18
$ws.Range("E18").Value = "  -0.65%  "

# Row 19
$ws.Range("E19").Value = "  -1.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.02"
$ws.Range("E20").Value = "  -6.09%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "444.06"
$ws.Range("E21").Value = "  -5.89%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.02"
$ws.Range("E22").Value = "  -13.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.623"
$ws.Range("E23").Value = "  -5.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.08"
$ws.Range("E24").Value = "  -3.53%  "

# Row 25
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.630.60"
$ws.Range("E26").Value = "  -5.28%  "

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000126"
$ws.Range("E27").Value = "  -1.56%  "

# Row 28
$ws.Range("E28").Value = "  -9.64%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.23"
$ws.Range("E29").Value = "  -7.89%  "

# Row 30
$ws.Range("E30").Value = "  -4.91%  "

# Row 31
$ws.Range("E31").Value = "  -7.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.09%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.164"
$ws.Range("E33").Value = "  -1.66%  "

# Row 34
$ws.Range("E34").Value = "  -4.51%  "

# Row 35
$ws.Range("E35").Value = "  -6.41%  "

# Row 36
$ws.Range("E36").Value = "  -7.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.474.84"
$ws.Range("E37").Value = "  -5.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.99"
$ws.Range("E38").Value = "  -5.22%  "

# Row 39
$ws.Range("E39").Value = "  +0.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.07%  "

# Row 41
$ws.Range("E41").Value = "  -0.52%  "

# Row 42
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0869"
$ws.Range("E42").Value = "  -3.75%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "168.38"
$ws.Range("E43").Value = "  -5.14%  "

# Row 44
$ws.Range("E44").Value = "  -8.08%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.883"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.66"
$ws.Range("E46").Value = "  -2.30%  "

# Row 47
$ws.Range("E47").Value = "  -0.89%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.21"
$ws.Range("E48").Value = "  -10.07%  "

# Row 49
$ws.Range("E49").Value = "  -11.18%  "

# Row 50
$ws.Range("E50").Value = "  -4.31%  "

# Row 51
$ws.Range("E51").Value = "  -4.49%  "
